$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3213")

# --- Scalar field updates (value changed, same label) ---
$ws.Range("B5").Value = "2"
$ws.Range("C5").Value = "2"

$ws.Range("B7").Value = "30 h"
$ws.Range("C7").Value = "30 h"

$ws.Range("B8").Value = "01/01/2023"
$ws.Range("C8").Value = "01/01/2023"

$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# --- New content added to row 11 (Objectives:) ---
$objectives = "Introduction of concepts related to rate and flow of momentum, heat and mass. Enable the student to model and solve problems of interest in transport phenomena, with appropriate choice of hypotheses and application of corresponding solution tools."
$ws.Range("B11").Value = $objectives
$ws.Range("C11").Value = $objectives

# --- New content added to row 13 (Programa resumido:) ---
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- New content added to row 14 (Short syllabus:) ---
$shortSyllabus = "Introduction to heat transfer. Heat conduction in steady state and in transient regime. Free and forced convection heat transfer. Heat transfer by thermal radiation. Heat transfer with phase change. Mass transference. Application examples."
$ws.Range("B14").Value = $shortSyllabus
$ws.Range("C14").Value = $shortSyllabus

# --- Row 15 (Programa:) reuses the "01/01/2023" string, per source data ---
$ws.Range("B15").Value = "01/01/2023"
$ws.Range("C15").Value = "01/01/2023"

# --- New content added to row 16 (Syllabus:) ---
$syllabus = "• Heat transfer by conduction: one-dimensional heat transfer in steady state. Fourier equation. Thermal conductivity. • One-dimensional heat transfer in steady state with convective contours. Newton's Law of Cooling. • Transient heat conduction. Thermal diffusivity. Biot number. • Analogy between heat transfer and electrical circuits: concepts of thermal resistance and capacitance. • Transient 2D and 3D heat transfer. • Convection heat transfer. Free convection. Similarity parameters. Rayleigh number. Forced convection. Boundary layer theory. Prandtl number and Nusselt number. • Radiation heat transfer. Blackbody radiation. Radiation properties. Radiation form factor. • Phase change heat transfer: boiling and condensation."
$ws.Range("B16").Value = $syllabus
$ws.Range("C16").Value = $syllabus

# --- Row 18 (Método:) reuses "1176388 - Luiz Tadeu Fernandes Eleno" ---
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
